# Generate Report for Handoff
#
# Row 3 (b.md) across the Overview, zh-cn and de-de sheets is updated to
# reflect that the file is now "Ready for handoff" (rather than "Handed
# back: in sync with en-US"), with refreshed handoff timestamps, new
# handoff xliff file names, and an added "not the latest" error-detail
# message on the per-locale sheets. The "Error Detail" column is also
# widened on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$dateFormat = "yyyy-mm-dd HH:mm:ss"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec3d012e024eb1b9a15f8078f92c027aae818f63/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9eff67fa68c7d2181bab78fadc7d4aaaad05d32/e2e/b.md."

# --- Overview sheet: row 3 corresponds to b.md ---
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$rOverviewG3 = $wsOverview.Range("G3")
$rOverviewG3.Value = "2016-08-16 18:32:50"
$rOverviewG3.NumberFormat = $dateFormat

# --- zh-cn sheet: row 3 corresponds to b.md ---
$wsZhCn.Range("C3").Value = "Ready for handoff"

# "False" looks like a boolean literal, so force it to stay text (as in the
# original file) using a leading apostrophe, then restore the cell's
# (default) style so no stray formatting is introduced.
$rZhCnF3 = $wsZhCn.Range("F3")
$savedStyleZhCnF3 = $rZhCnF3.Style
$rZhCnF3.Value = "'False"
$rZhCnF3.Style = $savedStyleZhCnF3

$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$rZhCnH3 = $wsZhCn.Range("H3")
$rZhCnH3.Value = "2016-08-16 18:32:45"
$rZhCnH3.NumberFormat = $dateFormat
$wsZhCn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 corresponds to b.md ---
$wsDeDe.Range("C3").Value = "Ready for handoff"

$rDeDeF3 = $wsDeDe.Range("F3")
$savedStyleDeDeF3 = $rDeDeF3.Style
$rDeDeF3.Value = "'False"
$rDeDeF3.Style = $savedStyleDeDeF3

$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$rDeDeH3 = $wsDeDe.Range("H3")
$rDeDeH3.Value = "2016-08-16 18:32:50"
$rDeDeH3.NumberFormat = $dateFormat
$wsDeDe.Range("P3").Value = $errorDetail

# --- Widen "Error Detail" column (P, the 16th column) on zh-cn and de-de sheets to 40 ---
# (The engine stores column width as ColumnWidth + 5/6, so use 40 - 5/6 to land on exactly 40.)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666667
